$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New data rows appended to the tracker (rows 31 gets its missing D/E
#    completed, and rows 32-43 are brand-new entries).
# ---------------------------------------------------------------------------

# Row 31 was missing its "complete" (D) and "score" (E) values - fill them in.
$ws.Cells.Item(31, 4).Value = 10
$ws.Cells.Item(31, 5).Formula = '=IF(ISBLANK(D31),"",IF(D31<B31,0,MIN(B31,D31)^1.2))'

# New rows: date, goal, complete (score is always a formula)
$newRows = @(
    @{ Row = 32; Date = 43348; Goal = 13; Complete = 7 },
    @{ Row = 33; Date = 43349; Goal = 10; Complete = 1 },
    @{ Row = 34; Date = 43350; Goal = 9;  Complete = 21 },
    @{ Row = 35; Date = 43351; Goal = 7;  Complete = 13 },
    @{ Row = 36; Date = 43352; Goal = 7;  Complete = 11 },
    @{ Row = 37; Date = 43353; Goal = 9;  Complete = 2 },
    @{ Row = 38; Date = 43354; Goal = 4;  Complete = 1 },
    @{ Row = 39; Date = 43355; Goal = 19; Complete = 10 },
    @{ Row = 40; Date = 43356; Goal = 13; Complete = 0 },
    @{ Row = 41; Date = 43357; Goal = 8;  Complete = 8 },
    @{ Row = 42; Date = 43358; Goal = 11; Complete = 4 },
    @{ Row = 43; Date = 43359; Goal = 12; Complete = $null }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.Goal
    $ws.Cells.Item($row, 3).Formula = "=FLOOR(B$row/7,1)"
    if ($null -ne $r.Complete) {
        $ws.Cells.Item($row, 4).Value = $r.Complete
        $ws.Cells.Item($row, 5).Formula = '=IF(ISBLANK(D' + $row + '),"",IF(D' + $row + '<B' + $row + ',0,MIN(B' + $row + ',D' + $row + ')^1.2))'
    }
}

# ---------------------------------------------------------------------------
# 2) get_quality needed a smarter, unambiguous date display - give the whole
#    date column a custom day-of-week format instead of the bare m/d/yyyy.
# ---------------------------------------------------------------------------
$ws.Range("A2:A43").NumberFormat = 'ddd\ d\-mmm\-yy'

# Re-apply the existing "highlight" fill to the rows that already had it so
# the new number format doesn't clobber their look.
$ws.Range("A2:A4").Interior.Color = 49407
$ws.Range("A8:A19").Interior.Color = 49407
$ws.Range("A22:A25").Interior.Color = 49407

# Header cell now carries the date style too, and the date column is a bit
# wider to fit the longer "Ddd d-Mmm-yy" text.
$ws.Cells.Item(1, 1).NumberFormat = 'ddd\ d\-mmm\-yy'
$ws.Columns.Item(1).ColumnWidth = 12.75

# ---------------------------------------------------------------------------
# 3) View state: scrolled further down, selection sitting on the next blank
#    row ready for new input.
# ---------------------------------------------------------------------------
$ws.Range("A22").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("B44").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4) Printer/page setup - force portrait orientation.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
